$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate column K (K4:K12) into column L (L4:L12), carrying over
# values and formatting, to add the "2020" column to the table.
$src = $ws.Range("K4:K12")
$dst = $ws.Range("L4:L12")
$src.Copy($dst)

# The new year value for the appended column.
$ws.Range("L4").Value = 2020

# Update the active selection as recorded after the edit.
$ws.Range("N5").Select()
